$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Hora (G) columns store plain text numbers/strings (e.g. "3.380",
# "12") as inline strings in the original workbook. Force text format first so
# Excel does not silently coerce the assigned strings into numeric values
# (which would strip significant trailing/leading zeros).
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Every data row (2-51): the scrape hour "Hora" moved from 12 to 13.
$ws.Range("G2:G51").Value = "13"

# Row 2
$ws.Range("D2").Value = '247.49'

# Row 4
$ws.Range("D4").Value = '5.525'

# Row 5
$ws.Range("D5").Value = '0.05638'

# Row 6
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").Value = '6.469'
$ws.Range("E6").Value = '5KuCoinTokenKCS'

# Row 7
$ws.Range("B7").Value = 'MXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D7").Value = '0.8036'
$ws.Range("E7").Value = '6MXTokenMX'

# Row 8
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").Value = '1.049'
$ws.Range("E8").Value = '7FTXTokenFTT'

# Row 9
$ws.Range("B9").Value = 'One'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D9").Value = '0.01171'
$ws.Range("E9").Value = '8OneONEBestin24h'

# Row 10
$ws.Range("D10").Value = '0.1424'

# Row 11
$ws.Range("D11").Value = '0.07324'

# Row 12
$ws.Range("D12").Value = '0.03191'

# Row 13
$ws.Range("D13").Value = '0.02947'

# Row 14
$ws.Range("D14").Value = '0.09265'

# Row 15
$ws.Range("D15").Value = '0.001662'

# Row 16
$ws.Range("D16").Value = '3.213'

# Row 17
$ws.Range("D17").Value = '0.04704'

# Row 18
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").Value = '0.006267'
$ws.Range("E18").Value = '17TigerCashTCH'

# Row 19
$ws.Range("B19").Value = 'BitKan'
$ws.Range("C19").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D19").Value = '0.001051'
$ws.Range("E19").Value = '18BitKanKAN'

# Row 20
$ws.Range("B20").Value = 'HotbitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D20").Value = '0.004121'
$ws.Range("E20").Value = '19HotbitTokenHTB'

# Row 21
$ws.Range("B21").Value = 'NitroEx'
$ws.Range("C21").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D21").Value = '0.0001501'
$ws.Range("E21").Value = '20NitroExNTX'

# Row 22
$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D22").Value = '3.973'
$ws.Range("E22").Value = '21LEOLEO'

# Row 23
$ws.Range("B23").Value = 'GateToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D23").Value = '3.380'
$ws.Range("E23").Value = '22GateTokenGT'

# Row 24
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").Value = '2.084'
$ws.Range("E24").Value = '23BTSETokenBTSE'

# Row 25
$ws.Range("B25").Value = 'BitpandaEcosystemToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D25").Value = '0.3272'
$ws.Range("E25").Value = '24BitpandaEcosystemTokenBEST'

# Row 26
$ws.Range("B26").Value = 'ProBitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D26").Value = '0.1285'
$ws.Range("E26").Value = '25ProBitTokenPROB'

# Row 27
$ws.Range("B27").Value = 'UpBots'
$ws.Range("C27").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D27").Value = '0.0003303'
$ws.Range("E27").Value = '26UpBotsUBXT'

# Row 40
$ws.Range("D40").Value = '0.04153'

# Row 41
$ws.Range("D41").Value = '0.006907'

# Row 42
$ws.Range("D42").Value = '0.003503'

# Row 43
$ws.Range("D43").Value = '0.1041'

# Row 44
$ws.Range("D44").Value = '0.009854'

# Row 45
$ws.Range("D45").Value = '0.00005645'

# Row 47
$ws.Range("D47").Value = '0.6807'

# Row 48
$ws.Range("D48").Value = '0.02465'
